$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the "last row" date-only style by moving it to the new row 39,
# and give row 38 the regular datetime style (same as the rest of the data rows).
$lastRowFormat = $ws.Cells.Item(38, 1).NumberFormat   # "YYYY-MM-DD" style (s=3)
$dataRowFormat = $ws.Cells.Item(37, 1).NumberFormat   # "YYYY-MM-DD HH:MM:SS" style (s=2)

$ws.Cells.Item(38, 1).NumberFormat = $dataRowFormat

$ws.Cells.Item(39, 1).Value = 45779
$ws.Cells.Item(39, 1).NumberFormat = $lastRowFormat
$ws.Cells.Item(39, 2).Value = 158
$ws.Cells.Item(39, 3).Value = 167
$ws.Cells.Item(39, 4).Value = 159
